# Generate Report for Archive
#
# The localization status moved from "Ready for handoff" to "In Translation"
# for the tracked file. That status string shows up on every sheet:
#   - Overview!E2 and Overview!F2 (the zh-cn / de-de roll-up columns)
#   - zh-cn!C2  (Status column)
#   - de-de!C2  (Status column)
#
# Because the new text is shorter than the old text, the Status columns
# (Overview E:F, and column C on the per-locale sheets) also shrink to
# fit the new content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
